# Script C ("Plot power output (single HPP)") gains a new switch,
# plot_RoR_part, that lets the user toggle the run-of-river component of
# the generation plot on (1) or off (0). This mirrors the existing
# single-value control rows already on the sheet (rows 1-5): column A
# holds the parameter name, column B is the (initially blank) user-input
# cell, and column C holds the explanatory text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plot power output (single HPP)")

# Fill the new row's text first (description in C, then the parameter
# name in A) so new shared-string entries are minted in the same order
# as they appear in the finished workbook.
$ws.Range("C6").Value = "set to 1 to plot run-of-river component of electricity generation, or to 0 to leave out this component"
$ws.Range("A6").Value = "plot_RoR_part"

# B6 stays empty (user fills it in later) but should carry the same
# "Note" cell style already used for the blank B1 input cell, rather
# than plain "Input" styling like B2:B5.
$ws.Range("B1").Copy()
$ws.Range("B6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Reflect where the user's selection ended up after making the edit.
$ws.Range("D8").Select()
